$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty row 23, which shifts rows 24-33 up by one.
# This moves "Feedback from users during test" (and the feedback list below it)
# up one row, closing the gap between row 22 and the feedback block.
$ws.Rows("23").Delete()

# Adjust column widths A:F to match the new layout (values taken from target file,
# pre-compensated for this engine's ColumnWidth -> stored-width conversion so the
# resulting OOXML <col width="..."/> lands as close as possible to the target)
$ws.Columns("A").ColumnWidth = 18.333333333333332
$ws.Columns("B").ColumnWidth = 29.5
$ws.Columns("C").ColumnWidth = 16.0
$ws.Columns("D").ColumnWidth = 29.166666666666668
$ws.Columns("E").ColumnWidth = 26.5
$ws.Columns("F").ColumnWidth = 24.166666666666668
